# Macroferia Regional de Talca - Limon : add a new day of price records.
#
# The workbook stores one flat table of daily price observations.  A new
# day's worth of data (3 records, date serial 44474 = 2021-10-05) is being
# inserted right after the existing block of rows for date 44273 (row 630),
# which pushes every following row down by 3 positions (old row N becomes
# new row N+3).  We replicate that by inserting 3 blank rows at row 631 and
# then filling them in with the new records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 631, shifting rows 631:723 down
# to 634:726 (dimension grows from A1:T723 to A1:T726).
$ws.Rows("631:633").Insert()

# Columns A-K are identical for every record in this sheet, so reuse the
# same constant values for the 3 new rows.
$category   = 5
$market     = "Macroferia Regional de Talca"
$region     = "Maule"
$newDate    = 44474
$weekday    = 7
$group      = "Fruta"
$famCode    = 100102
$family     = "Cítricos"
$prodCode   = 100102003
$product    = "Limón"
$variety    = "Sin especificar"

$newRows = @(
    @{ Row=631; L="1a amarillo"; M=360; N=5000; O=5000; P=5000; Q="`$/malla 16 kilos"; R="Cabildo";                  S=312; T=16 },
    @{ Row=632; L="1a amarillo"; M=200; N=6000; O=6000; P=6000; Q="`$/malla 16 kilos"; R="Provincia de Limarí";      S=375; T=16 },
    @{ Row=633; L="1a amarillo"; M=450; N=6000; O=6000; P=6000; Q="`$/malla 16 kilos"; R="Provincia de Melipilla";   S=375; T=16 }
)

foreach ($rec in $newRows) {
    $r = $rec.Row

    $ws.Cells.Item($r, 1).Value  = $category
    $ws.Cells.Item($r, 2).Value  = $market
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $newDate
    $ws.Cells.Item($r, 5).Value  = $weekday
    $ws.Cells.Item($r, 6).Value  = $group
    $ws.Cells.Item($r, 7).Value  = $famCode
    $ws.Cells.Item($r, 8).Value  = $family
    $ws.Cells.Item($r, 9).Value  = $prodCode
    $ws.Cells.Item($r, 10).Value = $product
    $ws.Cells.Item($r, 11).Value = $variety
    $ws.Cells.Item($r, 12).Value = $rec.L
    $ws.Cells.Item($r, 13).Value = $rec.M
    $ws.Cells.Item($r, 14).Value = $rec.N
    $ws.Cells.Item($r, 15).Value = $rec.O
    $ws.Cells.Item($r, 16).Value = $rec.P
    $ws.Cells.Item($r, 17).Value = $rec.Q
    $ws.Cells.Item($r, 18).Value = $rec.R
    $ws.Cells.Item($r, 19).Value = $rec.S
    $ws.Cells.Item($r, 20).Value = $rec.T
}
